# Atualização automática: 2025-09-05 09:00:25
# Rows 7-11 of the detections table got re-ordered: the record that used to be
# on row 11 (Placa PLACA_20250723145134 / Moura) moved up to row 7, and the
# four "Beja" records that used to occupy rows 7-10 each shifted down by one
# row (7->8, 8->9, 9->10, 10->11). Columns B (Class) and C (First_Detection_Date)
# are identical across these rows, so only A, D, E, F, G, H, I, J need updating.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Row, $A, $D, $E, $F, $G, $H, $I, $J)

    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H

    # Columns I ("First_Coords", e.g. "702,633,740,690") and J
    # ("First_Confidence", e.g. "0.76") are stored as text, not numbers.
    # Some of these comma/decimal separated values look numeric to Excel's
    # automatic type detection (e.g. "702,633,740,690" parses as a
    # thousands-grouped number). Force text formatting before assignment so
    # Excel does not coerce the value into a numeric cell, then restore the
    # default style so the cell formatting matches the rest of the sheet.
    $iCell = $ws.Cells.Item($Row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value = $I
    $iCell.Style = "Normal"

    $jCell = $ws.Cells.Item($Row, 10)
    $jCell.NumberFormat = "@"
    $jCell.Value = $J
    $jCell.Style = "Normal"
}

Set-RowValues 7  "2117575c-4ae1-458c-b88a-fc40f40debdb" "image_20250727074723_ppp0.jpg" "PLACA_20250723145134" "Moura" 38.06587 -7.221796 "1490,161,1563,258" "0.62"
Set-RowValues 8  "283b6eda-9c83-4cdd-9524-c7c394f2dc89" "image_20250728214139_ppp0.jpg" "PLACA_20250717165933" "Beja"  38.02035 -7.94715   "962,713,1006,765"  "0.76"
Set-RowValues 9  "a19b65d1-6f97-4841-9e1c-7446a9be92b6" "image_20250728214139_ppp0.jpg" "PLACA_20250717165933" "Beja"  38.02035 -7.94715   "967,614,1002,659"  "0.73"
Set-RowValues 10 "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d" "image_20250728214139_ppp0.jpg" "PLACA_20250717165933" "Beja"  38.02035 -7.94715   "702,633,740,690"   "0.72"
Set-RowValues 11 "dfd476d4-7689-4671-a076-78fe3ce806bb" "image_20250728214139_ppp0.jpg" "PLACA_20250717165933" "Beja"  38.02035 -7.94715   "1254,850,1294,895" "0.67"
